$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill rows 1176-1188 with new session entries (Sessions sheet) ---
$ws.Range("A1176").Value = "2024-08-19"
$ws.Range("B1176").Value = "09:45"
$ws.Range("C1176").Value = "10:45"
$ws.Range("D1176").Value = "1h 00m"
$ws.Range("E1176").Value = "#maintenance"
$ws.Range("G1176").Value = "'False"
$ws.Range("H1176").Value = "'False"
$ws.Range("I1176").Formula = "=YEAR(A1176)"
$ws.Range("J1176").Formula = "=MONTH(A1176)"

$ws.Range("A1177").Value = "2024-08-23"
$ws.Range("B1177").Value = "08:15"
$ws.Range("C1177").Value = "08:45"
$ws.Range("D1177").Value = "0h 30m"
$ws.Range("E1177").Value = "#python"
$ws.Range("F1177").Value = "nwragmate v1.0.0"
$ws.Range("G1177").Value = "'True"
$ws.Range("H1177").Value = "'False"
$ws.Range("I1177").Formula = "=YEAR(A1177)"
$ws.Range("J1177").Formula = "=MONTH(A1177)"

$ws.Range("A1178").Value = "2024-08-23"
$ws.Range("B1178").Value = "17:00"
$ws.Range("C1178").Value = "17:30"
$ws.Range("D1178").Value = "0h 30m"
$ws.Range("E1178").Value = "#python"
$ws.Range("F1178").Value = "nwragmate v1.0.0"
$ws.Range("G1178").Value = "'True"
$ws.Range("H1178").Value = "'False"
$ws.Range("I1178").Formula = "=YEAR(A1178)"
$ws.Range("J1178").Formula = "=MONTH(A1178)"

$ws.Range("A1179").Value = "2024-08-24"
$ws.Range("B1179").Value = "15:00"
$ws.Range("C1179").Value = "19:00"
$ws.Range("D1179").Value = "4h 00m"
$ws.Range("E1179").Value = "#maintenance"
$ws.Range("G1179").Value = "'False"
$ws.Range("H1179").Value = "'False"
$ws.Range("I1179").Formula = "=YEAR(A1179)"
$ws.Range("J1179").Formula = "=MONTH(A1179)"

$ws.Range("A1180").Value = "2024-08-31"
$ws.Range("B1180").Value = "16:00"
$ws.Range("C1180").Value = "18:00"
$ws.Range("D1180").Value = "2h 00m"
$ws.Range("E1180").Value = "#python"
$ws.Range("F1180").Value = "pycaretlab v1.0.0"
$ws.Range("G1180").Value = "'True"
$ws.Range("H1180").Value = "'False"
$ws.Range("I1180").Formula = "=YEAR(A1180)"
$ws.Range("J1180").Formula = "=MONTH(A1180)"

$ws.Range("A1181").Value = "2024-09-02"
$ws.Range("B1181").Value = "15:00"
$ws.Range("C1181").Value = "18:00"
$ws.Range("D1181").Value = "3h 00m"
$ws.Range("E1181").Value = "#python"
$ws.Range("F1181").Value = "nwragmate v1.0.0"
$ws.Range("G1181").Value = "'True"
$ws.Range("H1181").Value = "'False"
$ws.Range("I1181").Formula = "=YEAR(A1181)"
$ws.Range("J1181").Formula = "=MONTH(A1181)"

$ws.Range("A1182").Value = "2024-09-02"
$ws.Range("B1182").Value = "21:00"
$ws.Range("C1182").Value = "22:45"
$ws.Range("D1182").Value = "1h 45m"
$ws.Range("E1182").Value = "#python"
$ws.Range("F1182").Value = "nwragmate v1.0.0"
$ws.Range("G1182").Value = "'True"
$ws.Range("H1182").Value = "'False"
$ws.Range("I1182").Formula = "=YEAR(A1182)"
$ws.Range("J1182").Formula = "=MONTH(A1182)"

$ws.Range("A1183").Value = "2024-09-03"
$ws.Range("B1183").Value = "09:30"
$ws.Range("C1183").Value = "13:00"
$ws.Range("D1183").Value = "3h 30m"
$ws.Range("E1183").Value = "#python"
$ws.Range("F1183").Value = "nwragmate v1.0.0"
$ws.Range("G1183").Value = "'True"
$ws.Range("H1183").Value = "'False"
$ws.Range("I1183").Formula = "=YEAR(A1183)"
$ws.Range("J1183").Formula = "=MONTH(A1183)"

$ws.Range("A1184").Value = "2024-09-03"
$ws.Range("B1184").Value = "19:30"
$ws.Range("C1184").Value = "22:00"
$ws.Range("D1184").Value = "2h 30m"
$ws.Range("E1184").Value = "#python"
$ws.Range("F1184").Value = "nwragmate v1.0.0"
$ws.Range("G1184").Value = "'True"
$ws.Range("H1184").Value = "'False"
$ws.Range("I1184").Formula = "=YEAR(A1184)"
$ws.Range("J1184").Formula = "=MONTH(A1184)"

$ws.Range("A1185").Value = "2024-09-05"
$ws.Range("B1185").Value = "08:00"
$ws.Range("C1185").Value = "08:30"
$ws.Range("D1185").Value = "0h 30m"
$ws.Range("E1185").Value = "#python"
$ws.Range("F1185").Value = "nwragmate v1.0.0"
$ws.Range("G1185").Value = "'True"
$ws.Range("H1185").Value = "'False"
$ws.Range("I1185").Formula = "=YEAR(A1185)"
$ws.Range("J1185").Formula = "=MONTH(A1185)"

$ws.Range("A1186").Value = "2024-09-05"
$ws.Range("B1186").Value = "17:00"
$ws.Range("C1186").Value = "17:30"
$ws.Range("D1186").Value = "0h 30m"
$ws.Range("E1186").Value = "#python"
$ws.Range("F1186").Value = "nwragmate v1.0.0"
$ws.Range("G1186").Value = "'True"
$ws.Range("H1186").Value = "'False"
$ws.Range("I1186").Formula = "=YEAR(A1186)"
$ws.Range("J1186").Formula = "=MONTH(A1186)"

$ws.Range("A1187").Value = "2024-09-06"
$ws.Range("B1187").Value = "08:00"
$ws.Range("C1187").Value = "08:45"
$ws.Range("D1187").Value = "0h 45m"
$ws.Range("E1187").Value = "#python"
$ws.Range("F1187").Value = "nwragmate v1.0.0"
$ws.Range("G1187").Value = "'True"
$ws.Range("H1187").Value = "'False"
$ws.Range("I1187").Formula = "=YEAR(A1187)"
$ws.Range("J1187").Formula = "=MONTH(A1187)"

$ws.Range("A1188").Value = "2024-09-06"
$ws.Range("B1188").Value = "17:00"
$ws.Range("C1188").Value = "17:30"
$ws.Range("D1188").Value = "0h 30m"
$ws.Range("E1188").Value = "#python"
$ws.Range("F1188").Value = "nwragmate v1.0.0"
$ws.Range("G1188").Value = "'True"
$ws.Range("H1188").Value = "'False"
$ws.Range("I1188").Formula = "=YEAR(A1188)"
$ws.Range("J1188").Formula = "=MONTH(A1188)"

# --- Append new trailing blank rows (1192-1202), matching the style of the
#     existing blank rows further up (e.g. row 1189) ---
$ws.Range("A1189:J1189").Copy()
$ws.Range("A1192:J1202").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update active selection to mirror the author's final view state
#     (keep the existing "freeze header row" pane layout intact) ---
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F1197").Select()
